$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing data (old A..S) to B..T.
$ws.Columns("A:A").Insert()

# Populate the new "Fund" filter column used for the 40017-B trades.
$ws.Range("A4").Value = "Fund"
$ws.Range("A5").Value = "40017-B"
$ws.Range("A6").Value = "40017-B"
$ws.Range("A7").Value = "40017-B"

# Match the cell formatting used by the rest of the data (copy format from col B).
$ws.Range("B4:B7").Copy()
$ws.Range("A4:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Columns("A:A").ColumnWidth = 8.67

# Update the active selection like in the edited workbook.
$ws.Range("A8").Select()
